$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8924640417098999
$ws.Range("B1").Value = 1.77788233757019
$ws.Range("C1").Value = 4.12475061416626
$ws.Range("D1").Value = 3.549455642700195
$ws.Range("E1").Value = 1.509981989860535
